$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (rows 2-7): add "Unique Subcodes" / "Modification Date" columns ---

# Header row (row 2): copy style from E2 (existing last header cell) onto the
# two new header cells and set their text.
$ws.Range("F2").Value = "Unique Subcodes"
$ws.Range("G2").Value = "Modification Date"
$ws.Range("E2").Copy()
$ws.Range("F2:G2").PasteSpecial(-4122)

# Data row (row 3): new boolean-looking "FALSE" text cell + new timestamp text
# cell. Both values must land as literal text (not auto-coerced types), so
# write them first as formulas on a scratch cell and paste only the computed
# *values* across - this keeps them as plain strings.
$scratch = $ws.Range("Z1")

$scratch.Formula = "=""FALSE"""
$scratch.Copy()
$ws.Range("F3").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("F3").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

$scratch.Formula = "=""2023-03-11 17:23:44"""
$scratch.Copy()
$ws.Range("G3").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("G3").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# --- Table 2 (rows 10-15): same two new columns ---

$ws.Range("F10").Value = "Unique Subcodes"
$ws.Range("G10").Value = "Modification Date"
$ws.Range("E10").Copy()
$ws.Range("F10:G10").PasteSpecial(-4122)

$scratch.Formula = "=""FALSE"""
$scratch.Copy()
$ws.Range("F11").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("F11").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

$scratch.Formula = "=""2023-03-11 17:23:44"""
$scratch.Copy()
$ws.Range("G11").PasteSpecial(-4163)
$scratch.Clear()
$ws.Range("G11").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

$scratch.Clear()

# Match the author's final selection.
$ws.Range("F3").Select()
